$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "resistenze composite parallelo" sheet (sheet5): drop the stray empty
#    A4 cell that carries no value.
# ---------------------------------------------------------------------------
$wsParallelo = $wb.Worksheets.Item("resistenze composite parallelo")
$wsParallelo.Range("A4").ClearContents()

# ---------------------------------------------------------------------------
# 2) Add the new "soglia diodo" sheet at the end of the workbook (after the
#    last existing sheet) and populate it with the collected diode-threshold
#    measurements.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "soglia diodo"

# Column widths
$ws.Columns.Item(3).ColumnWidth = 46

# Header row
$ws.Range("A1").Value = "V"
$ws.Range("B1").Value = "I(uA)"
$ws.Range("C1").Value = "NOTA: attorno ai 0.7-0.8 i valori sono sulla soglia"

# Data rows: col A = V reading (3 decimals), col B = I(uA) reading
$data = @(
    @(2,  0.06,   0.001),
    @(3,  0.106,  0.005),
    @(4,  0.152,  0.005),
    @(5,  0.208,  0.01),
    @(6,  0.249,  0.025),
    @(7,  0.3,    0.105),
    @(8,  0.356,  0.592),
    @(9,  0.403,  2.54),
    @(10, 0.454,  11.89),
    @(11, 0.5,    40.27),
    @(12, 0.556,  123.63),
    @(13, 0.602,  234.85),
    @(14, 0.649,  375.42),
    @(15, 0.704,  574.25),
    @(16, 0.745,  732.1),
    @(17, 0.797,  937.7),
    @(18, 0.845,  36863),
    @(19, 0.899,  53710),
    @(20, 0.954,  71790),
    @(21, 0.999,  87254)
)

foreach ($row in $data) {
    $r = $row[0]
    $vVal = $row[1]
    $iVal = $row[2]

    $aCell = $ws.Range("A$r")
    $aCell.Value = $vVal
    $aCell.NumberFormat = "0.000"

    $bCell = $ws.Range("B$r")
    $bCell.Value = $iVal
    if ($r -ge 18) {
        $bCell.NumberFormat = "0"
    } elseif ($r -ge 12) {
        $bCell.NumberFormat = "0.00"
    } else {
        $bCell.NumberFormat = "0.000"
    }
}
